# "Basics of arrays + fixes"
# Append three new bytecode instructions (newarray, aload, astore) to the
# VM opcode reference table on Sheet1, directly below the existing last
# row (row 58 / opcode 0x38 "syscall").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("0x39", "newarray", "Push a new array object to the stack"),
    @("0x3A", "aload",    "Load an object from the array ontop of the stack to the stack"),
    @("0x3B", "astore",   "Store an object ontop of the stack to and array now ontop of the stack")
)

$r = 59
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $r++
}

# Leave the selection on the last cell entered, matching the authoring
# session's final cursor position.
$ws.Range("D61").Select()
